$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the credential row for akhachatryan19881@mail.ru / arman1988 (row 8) by
# deleting the A8:C8 cell range and shifting the remaining cells up (rather than a
# full row delete), which matches the original author's edit: hyperlink refs for
# rows A2-A8 stay bound to their original relationship ids, and only the last
# (now-empty) hyperlink entry drops off.
$ws.Range("A8:C8").Delete(-4162)
